$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.946.60'
$ws.Range('E2').Value = '  +4.46%  '
$ws.Range('D3').Value = '1.879.61'
$ws.Range('E3').Value = '  +3.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '278.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5300'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3469'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.15'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06968'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.21'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.8093'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07859'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = '1.874.47'
$ws.Range('E14').Value = '  +3.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.199'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.20%  '
$ws.Range('E17').Value = '  +3.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9994'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008108'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = '26.969.41'
$ws.Range('E21').Value = '  +4.32%  '
$ws.Range('D22').Value = '2.100.87'
$ws.Range('E22').Value = '  +2.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.762'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.197'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.363'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '146.76'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.97%  '
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.663'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '113.92'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.376'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.342'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.08916'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04952'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.185'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7398'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.884'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.291'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.411'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01858'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5170'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9644'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '116.70'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.205'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.143'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9991'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4532'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1351'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.409'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.40%  '
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('E51').Value = '  +0.40%  '
